$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two factual prompts from Katy Perry / Uptown Girl to Taylor Swift / All Too Well
$ws.Range("C2").Value = "What do we know about Taylor Swift?"
$ws.Range("C3").Value = "Who sings All Too Well?"

# Column C had picked up a stray column-level format (an empty "applyAlignment"
# style with no actual alignment) that isn't used by any cell. Clear the whole
# column's formatting, then restore the real per-cell format (vertical
# center, the style actually used by the data cells) on the used range so
# only the orphaned column-level default is dropped.
$ws.Columns("C").ClearFormats()
$ws.Range("C1:C13").VerticalAlignment = -4108

$ws.Range("C5").Select()
